$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Table layout:
#   Row 1 = header (Distance Class (m) | N | Mantel r | p)
#   Row 2 = Distance Class 5,000
#   Row 3 = Distance Class 15,000
#   Row 4 = Distance Class 25,000
# Columns: 1 = Distance Class, 2 = N, 3 = Mantel r, 4 = p

# Row 2 (5,000): Mantel r -0.026 -> -0.017 ; p 0.135 -> 0.243
$table.Cell(2, 3).Range.Text = "-0.017"
$table.Cell(2, 4).Range.Text = "0.243"

# Row 3 (15,000): Mantel r -0.042 -> -0.026 ; p 0.056 -> 0.2
$table.Cell(3, 3).Range.Text = "-0.026"
$table.Cell(3, 4).Range.Text = "0.2"

# Row 4 (25,000): Mantel r 0.064 -> 0.007 ; p 0.021 -> 0.486 (no longer bold)
$table.Cell(4, 3).Range.Text = "0.007"

$pCell = $table.Cell(4, 4)
$pCell.Range.Text = "0.486"
$pStart = $pCell.Range.Start
$d.Range($pStart, $pStart + 5).Font.Bold = $false
